$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting (incl. date number format) from the row above, then fill values
$ws.Range("A13:H13").Copy($ws.Range("A14:H14"))

$ws.Range("A14").Value = 9392.8799999999992
$ws.Range("B14").Value = 9292.52
$ws.Range("C14").Value = 281.06
$ws.Range("D14").Value = 284.08999999999997
$ws.Range("E14").Value = $false
$ws.Range("F14").Value = 1.08
$ws.Range("G14").Value = 42620.766435185185
$ws.Range("H14").Value = $true
